# Generate Report for Handback
# The a9227bf9-... file has moved from "Ready for handoff" to
# "Handed back: in sync with en-US" for both the zh-cn and de-de locales.
# Update the Overview sheet plus the per-locale detail sheets accordingly,
# and record the new handback timestamps.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

# --- Overview sheet: File Name | zh-cn | de-de | Latest Handoff Date ---
$overview = $wb.Worksheets.Item("Overview")
# Row 3 corresponds to a9227bf9-d672-45db-85a6-1fe97592d078.md
$overview.Range("B3").Value = $statusHandedBack
$overview.Range("C3").Value = $statusHandedBack

# --- zh-cn sheet ---
$zhcn = $wb.Worksheets.Item("zh-cn")
# Row 3 corresponds to a9227bf9-d672-45db-85a6-1fe97592d078.md
$zhcn.Range("C3").Value = $statusHandedBack
$zhcn.Range("H3").Value = "2016-03-20 04:38:51"

# --- de-de sheet ---
$dede = $wb.Worksheets.Item("de-de")
# Row 3 corresponds to a9227bf9-d672-45db-85a6-1fe97592d078.md
$dede.Range("C3").Value = $statusHandedBack
$dede.Range("H3").Value = "2016-03-20 04:38:56"
